$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.804
$ws.Range("L2").Value = 0.499

$ws.Range("D3").Value = 0.493
$ws.Range("H3").Value = 0.375

$ws.Range("G4").Value = 0.734

$ws.Range("G5").Value = 0.789
$ws.Range("K5").Value = 0.417

$ws.Range("E6").Value = 0.598

$ws.Range("G62").Value = 0.578

$ws.Range("E65").Value = 0.603

$ws.Range("G81").Value = 0.565
$ws.Range("K81").Value = 0.382

$ws.Range("G84").Value = 0.771

$ws.Range("E97").Value = 0.519
$ws.Range("K97").Value = 0.358
